{"js": "// Replace each three-digit-by-one-digit division answer in the practice\n// table with a freshly generated problem (same \"a\u00f7b=c, r\" format).\nconst pairs = [\n  [\"176\u00f74=44, 0\", \"694\u00f78=86, 6\"],\n  [\"538\u00f77=76, 6\", \"269\u00f77=38, 3\"],\n  [\"566\u00f74=141, 2\", \"537\u00f76=89, 3\"],\n  [\"429\u00f77=61, 2\", \"930\u00f78=116, 2\"],\n  [\"221\u00f79=24, 5\", \"900\u00f75=180, 0\"],\n  [\"182\u00f76=30, 2\", \"872\u00f74=218, 0\"],\n  [\"797\u00f76=132, 5\", \"672\u00f72=336, 0\"],\n  [\"945\u00f73=315, 0\", \"453\u00f73=151, 0\"],\n  [\"589\u00f76=98, 1\", \"629\u00f76=104, 5\"],\n  [\"564\u00f74=141, 0\", \"959\u00f76=159, 5\"],\n  [\"818\u00f79=90, 8\", \"640\u00f75=128, 0\"],\n  [\"993\u00f77=141, 6\", \"445\u00f74=111, 1\"],\n  [\"577\u00f73=192, 1\", \"879\u00f74=219, 3\"],\n  [\"138\u00f79=15, 3\", \"215\u00f76=35, 5\"],\n  [\"736\u00f78=92, 0\", \"503\u00f73=167, 2\"],\n  [\"741\u00f79=82, 3\", \"482\u00f78=60, 2\"],\n  [\"879\u00f79=97, 6\", \"167\u00f72=83, 1\"],\n  [\"596\u00f75=119, 1\", \"440\u00f79=48, 8\"],\n  [\"318\u00f73=106, 0\", \"191\u00f78=23, 7\"],\n  [\"698\u00f75=139, 3\", \"521\u00f76=86, 5\"],\n  [\"439\u00f75=87, 4\", \"208\u00f77=29, 5\"],\n  [\"524\u00f77=74, 6\", \"270\u00f74=67, 2\"],\n  [\"311\u00f77=44, 3\", \"109\u00f78=13, 5\"],\n  [\"888\u00f74=222, 0\", \"292\u00f77=41, 5\"],\n  [\"631\u00f76=105, 1\", \"578\u00f78=72, 2\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"176\u00f74=44, 0\", \"694\u00f78=86, 6\"),\n    @(\"538\u00f77=76, 6\", \"269\u00f77=38, 3\"),\n    @(\"566\u00f74=141, 2\", \"537\u00f76=89, 3\"),\n    @(\"429\u00f77=61, 2\", \"930\u00f78=116, 2\"),\n    @(\"221\u00f79=24, 5\", \"900\u00f75=180, 0\"),\n    @(\"182\u00f76=30, 2\", \"872\u00f74=218, 0\"),\n    @(\"797\u00f76=132, 5\", \"672\u00f72=336, 0\"),\n    @(\"945\u00f73=315, 0\", \"453\u00f73=151, 0\"),\n    @(\"589\u00f76=98, 1\", \"629\u00f76=104, 5\"),\n    @(\"564\u00f74=141, 0\", \"959\u00f76=159, 5\"),\n    @(\"818\u00f79=90, 8\", \"640\u00f75=128, 0\"),\n    @(\"993\u00f77=141, 6\", \"445\u00f74=111, 1\"),\n    @(\"577\u00f73=192, 1\", \"879\u00f74=219, 3\"),\n    @(\"138\u00f79=15, 3\", \"215\u00f76=35, 5\"),\n    @(\"736\u00f78=92, 0\", \"503\u00f73=167, 2\"),\n    @(\"741\u00f79=82, 3\", \"482\u00f78=60, 2\"),\n    @(\"879\u00f79=97, 6\", \"167\u00f72=83, 1\"),\n    @(\"596\u00f75=119, 1\", \"440\u00f79=48, 8\"),\n    @(\"318\u00f73=106, 0\", \"191\u00f78=23, 7\"),\n    @(\"698\u00f75=139, 3\", \"521\u00f76=86, 5\"),\n    @(\"439\u00f75=87, 4\", \"208\u00f77=29, 5\"),\n    @(\"524\u00f77=74, 6\", \"270\u00f74=67, 2\"),\n    @(\"311\u00f77=44, 3\", \"109\u00f78=13, 5\"),\n    @(\"888\u00f74=222, 0\", \"292\u00f77=41, 5\"),\n    @(\"631\u00f76=105, 1\", \"578\u00f78=72, 2\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
